# Parser fonctionel? je crois bien
#
# The "Follow" set column (C) for the rows ExprArith, ExprArith', MultExpr,
# MultExpr' and Term gains a trailing ", do" token (the grammar now treats
# `do` as a valid Follow token for these productions).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value  = "…,else, end, ), *, /,  =, <, and, or, then, do"
$ws.Range("C9").Value  = "…,else, end ), *, /,  =, <, and, or, then, do"
$ws.Range("C10").Value = "…,else, end ), *, /,  =, <, and,or, then +, -, do"
$ws.Range("C11").Value = "…,else, end ), *, /,  =, <, and, or, then, +, -, do"
$ws.Range("C12").Value = "…,else, end ), *, /,  =, <, and, or, then, +, -, do"

# Move the active selection to C14 (was the whole A1:C26 range before).
$ws.Range("C14").Select()
